# Update "想去人数" (want-to-go count) values in column F for the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets, to
# reflect the refreshed numbers output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1612
$ws1.Range("F3").Value  = 8997
$ws1.Range("F4").Value  = 105
$ws1.Range("F6").Value  = 689
$ws1.Range("F7").Value  = 346
$ws1.Range("F8").Value  = 175
$ws1.Range("F10").Value = 79
$ws1.Range("F11").Value = 3837
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 378
$ws1.Range("F15").Value = 4275
$ws1.Range("F17").Value = 158
$ws1.Range("F18").Value = 1142
$ws1.Range("F19").Value = 9
$ws1.Range("F21").Value = 4
$ws1.Range("F22").Value = 241
$ws1.Range("F24").Value = 2643
$ws1.Range("F25").Value = 112

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 38

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1612
$ws4.Range("F3").Value  = 8997
$ws4.Range("F4").Value  = 105
$ws4.Range("F6").Value  = 689
$ws4.Range("F7").Value  = 346
$ws4.Range("F8").Value  = 175
$ws4.Range("F10").Value = 79
$ws4.Range("F11").Value = 3837
$ws4.Range("F12").Value = 57
$ws4.Range("F13").Value = 378
$ws4.Range("F15").Value = 4275
$ws4.Range("F17").Value = 158
$ws4.Range("F18").Value = 1143
$ws4.Range("F19").Value = 9
$ws4.Range("F21").Value = 4
$ws4.Range("F22").Value = 241
$ws4.Range("F24").Value = 2643
$ws4.Range("F25").Value = 38
$ws4.Range("F26").Value = 112
